$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7667979598045349
$ws.Range("B1").Value = 1.601015567779541
$ws.Range("C1").Value = 4.887022495269775
$ws.Range("D1").Value = 2.274394512176514
$ws.Range("E1").Value = 1.238748669624329
